$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 45090
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100114007
$ws.Range("G27").Value = "Jengibre"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 15000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 15600
$ws.Range("N27").Value = "`$/caja 13 kilos"
$ws.Range("O27").Value = "Perú"
$ws.Range("P27").Value = 1200
$ws.Range("Q27").Value = 13
$ws.Range("R27").Value = "Hortaliza"
